$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("E2").Value = 0.1010988340531174
$ws.Range("G2").Value = 0.0668636523347244
$ws.Range("H2").Value = 0.06564243869397875
$ws.Range("L2").Value = 0.01174687881814288
$ws.Range("M2").Value = 0.1387155700046751
$ws.Range("O2").Value = 0.1480593999244309
$ws.Range("R2").Value = 0.1296475057981424
$ws.Range("T2").Value = 0.2531094073799869
$ws.Range("V2").Value = 0.03208782705695191
$ws.Range("AC2").Value = 0.03858595623164139
$ws.Range("AE2").Value = 0.01444252970420799
$ws.Range("D3").Value = 0.09440458670469844
$ws.Range("F3").Value = 0.08928568358160086
$ws.Range("G3").Value = 0.1231695316628098
$ws.Range("H3").Value = 0.001016097278780906
$ws.Range("L3").Value = 0.06167035325117064
$ws.Range("N3").Value = 0.1854554872073387
$ws.Range("Q3").Value = 0.09609757495256303
$ws.Range("S3").Value = 0.2004404370403535
$ws.Range("T3").Value = 0.00400448844158176
$ws.Range("U3").Value = 0.04156828683272107
$ws.Range("X3").Value = 0.003799006799241321
$ws.Range("AB3").Value = 0.06398104488107292
$ws.Range("AD3").Value = 0.03510742136606692
$ws.Range("E4").Value = 0.1955346710943545
$ws.Range("F4").Value = 0.02503871942250605
$ws.Range("G4").Value = 0.10073926783001
$ws.Range("H4").Value = 0.07008884925579976
$ws.Range("M4").Value = 0.0731494926343001
$ws.Range("N4").Value = 0.09681703560635958
$ws.Range("O4").Value = 0.06343881840481565
$ws.Range("Q4").Value = 0.03530148627196269
$ws.Range("R4").Value = 0.07022628979129746
$ws.Range("S4").Value = 0.05215413099508406
$ws.Range("T4").Value = 0.1468558975630559
$ws.Range("U4").Value = 0.003397272532168043
$ws.Range("Y4").Value = 0.004297036873502315
$ws.Range("AB4").Value = 0.006172193742441649
$ws.Range("AC4").Value = 0.03648801970699914
$ws.Range("AD4").Value = 0.006998002194418473
$ws.Range("AE4").Value = 0.01330281608092464
$ws.Range("E5").Value = 0.1610451991759126
$ws.Range("G5").Value = 0.05234145078162396
$ws.Range("H5").Value = 0.07908271028594668
$ws.Range("L5").Value = 0.00122233208426858
$ws.Range("M5").Value = 0.1098990033659103
$ws.Range("N5").Value = 0.01724826651449359
$ws.Range("O5").Value = 0.1493301293472525
$ws.Range("R5").Value = 0.108745349979214
$ws.Range("T5").Value = 0.2091045593938385
$ws.Range("V5").Value = 0.01129712323415433
$ws.Range("Y5").Value = 0.008203695880124847
$ws.Range("AC5").Value = 0.07516797592595637
$ws.Range("AE5").Value = 0.01731220403130365
$ws.Range("E6").Value = 0.01215440627324932
$ws.Range("F6").Value = 0.04588114994516361
$ws.Range("G6").Value = 0.06080645156699901
$ws.Range("H6").Value = 0.1305609256222446
$ws.Range("I6").Value = 0.03482846621177704
$ws.Range("J6").Value = 0.01266075950643073
$ws.Range("M6").Value = 0.002374518522432611
$ws.Range("O6").Value = 0.2237667863129311
$ws.Range("R6").Value = 0.1278754036446064
$ws.Range("S6").Value = 0.008048243318566534
$ws.Range("T6").Value = 0.1348666170706091
$ws.Range("U6").Value = 0.03415564303509747
$ws.Range("V6").Value = 0.02949805686565234
$ws.Range("Z6").Value = 0.01178650308876554
$ws.Range("AC6").Value = 0.08223347141206504
$ws.Range("AD6").Value = 0.02004150957250685
$ws.Range("AE6").Value = 0.02846108803090245

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("E2").Value = 0.1010988340531174
$ws.Range("F2").Value = 0.1010988340531174
$ws.Range("G2").Value = 0.1679624863878418
$ws.Range("H2").Value = 0.2336049250818206
$ws.Range("I2").Value = 0.2336049250818206
$ws.Range("J2").Value = 0.2336049250818206
$ws.Range("K2").Value = 0.2336049250818206
$ws.Range("L2").Value = 0.2453518038999634
$ws.Range("M2").Value = 0.3840673739046385
$ws.Range("N2").Value = 0.3840673739046385
$ws.Range("O2").Value = 0.5321267738290694
$ws.Range("P2").Value = 0.5321267738290694
$ws.Range("Q2").Value = 0.5321267738290694
$ws.Range("R2").Value = 0.6617742796272118
$ws.Range("S2").Value = 0.6617742796272118
$ws.Range("T2").Value = 0.9148836870071987
$ws.Range("U2").Value = 0.9148836870071987
$ws.Range("V2").Value = 0.9469715140641506
$ws.Range("W2").Value = 0.9469715140641506
$ws.Range("X2").Value = 0.9469715140641506
$ws.Range("Y2").Value = 0.9469715140641506
$ws.Range("Z2").Value = 0.9469715140641506
$ws.Range("AA2").Value = 0.9469715140641506
$ws.Range("AB2").Value = 0.9469715140641506
$ws.Range("AC2").Value = 0.985557470295792
$ws.Range("AD2").Value = 0.985557470295792
$ws.Range("AE2").Value = 1
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 1
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 1
$ws.Range("AJ2").Value = 1
$ws.Range("D3").Value = 0.09440458670469844
$ws.Range("E3").Value = 0.09440458670469844
$ws.Range("F3").Value = 0.1836902702862993
$ws.Range("G3").Value = 0.3068598019491091
$ws.Range("H3").Value = 0.30787589922789
$ws.Range("I3").Value = 0.30787589922789
$ws.Range("J3").Value = 0.30787589922789
$ws.Range("K3").Value = 0.30787589922789
$ws.Range("L3").Value = 0.3695462524790606
$ws.Range("M3").Value = 0.3695462524790606
$ws.Range("N3").Value = 0.5550017396863993
$ws.Range("O3").Value = 0.5550017396863993
$ws.Range("P3").Value = 0.5550017396863993
$ws.Range("Q3").Value = 0.6510993146389623
$ws.Range("R3").Value = 0.6510993146389623
$ws.Range("S3").Value = 0.8515397516793158
$ws.Range("T3").Value = 0.8555442401208976
$ws.Range("U3").Value = 0.8971125269536186
$ws.Range("V3").Value = 0.8971125269536186
$ws.Range("W3").Value = 0.8971125269536186
$ws.Range("X3").Value = 0.9009115337528599
$ws.Range("Y3").Value = 0.9009115337528599
$ws.Range("Z3").Value = 0.9009115337528599
$ws.Range("AA3").Value = 0.9009115337528599
$ws.Range("AB3").Value = 0.9648925786339329
$ws.Range("AC3").Value = 0.9648925786339329
$ws.Range("AD3").Value = 0.9999999999999998
$ws.Range("AE3").Value = 0.9999999999999998
$ws.Range("AF3").Value = 0.9999999999999998
$ws.Range("AG3").Value = 0.9999999999999998
$ws.Range("AH3").Value = 0.9999999999999998
$ws.Range("AI3").Value = 0.9999999999999998
$ws.Range("AJ3").Value = 0.9999999999999998
$ws.Range("E4").Value = 0.1955346710943545
$ws.Range("F4").Value = 0.2205733905168606
$ws.Range("G4").Value = 0.3213126583468706
$ws.Range("H4").Value = 0.3914015076026704
$ws.Range("I4").Value = 0.3914015076026704
$ws.Range("J4").Value = 0.3914015076026704
$ws.Range("K4").Value = 0.3914015076026704
$ws.Range("L4").Value = 0.3914015076026704
$ws.Range("M4").Value = 0.4645510002369705
$ws.Range("N4").Value = 0.5613680358433301
$ws.Range("O4").Value = 0.6248068542481457
$ws.Range("P4").Value = 0.6248068542481457
$ws.Range("Q4").Value = 0.6601083405201085
$ws.Range("R4").Value = 0.7303346303114059
$ws.Range("S4").Value = 0.7824887613064899
$ws.Range("T4").Value = 0.9293446588695458
$ws.Range("U4").Value = 0.9327419314017138
$ws.Range("V4").Value = 0.9327419314017138
$ws.Range("W4").Value = 0.9327419314017138
$ws.Range("X4").Value = 0.9327419314017138
$ws.Range("Y4").Value = 0.9370389682752162
$ws.Range("Z4").Value = 0.9370389682752162
$ws.Range("AA4").Value = 0.9370389682752162
$ws.Range("AB4").Value = 0.9432111620176578
$ws.Range("AC4").Value = 0.9796991817246569
$ws.Range("AD4").Value = 0.9866971839190753
$ws.Range("E5").Value = 0.1610451991759126
$ws.Range("F5").Value = 0.1610451991759126
$ws.Range("G5").Value = 0.2133866499575366
$ws.Range("H5").Value = 0.2924693602434832
$ws.Range("I5").Value = 0.2924693602434832
$ws.Range("J5").Value = 0.2924693602434832
$ws.Range("K5").Value = 0.2924693602434832
$ws.Range("L5").Value = 0.2936916923277518
$ws.Range("M5").Value = 0.4035906956936622
$ws.Range("N5").Value = 0.4208389622081558
$ws.Range("O5").Value = 0.5701690915554083
$ws.Range("P5").Value = 0.5701690915554083
$ws.Range("Q5").Value = 0.5701690915554083
$ws.Range("R5").Value = 0.6789144415346222
$ws.Range("S5").Value = 0.6789144415346222
$ws.Range("T5").Value = 0.8880190009284608
$ws.Range("U5").Value = 0.8880190009284608
$ws.Range("V5").Value = 0.8993161241626151
$ws.Range("W5").Value = 0.8993161241626151
$ws.Range("X5").Value = 0.8993161241626151
$ws.Range("Y5").Value = 0.9075198200427399
$ws.Range("Z5").Value = 0.9075198200427399
$ws.Range("AA5").Value = 0.9075198200427399
$ws.Range("AB5").Value = 0.9075198200427399
$ws.Range("AC5").Value = 0.9826877959686963
$ws.Range("AD5").Value = 0.9826877959686963
$ws.Range("AE5").Value = 0.9999999999999999
$ws.Range("AF5").Value = 0.9999999999999999
$ws.Range("AG5").Value = 0.9999999999999999
$ws.Range("AH5").Value = 0.9999999999999999
$ws.Range("AI5").Value = 0.9999999999999999
$ws.Range("AJ5").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.01215440627324932
$ws.Range("F6").Value = 0.05803555621841293
$ws.Range("G6").Value = 0.1188420077854119
$ws.Range("H6").Value = 0.2494029334076566
$ws.Range("I6").Value = 0.2842313996194336
$ws.Range("J6").Value = 0.2968921591258644
$ws.Range("K6").Value = 0.2968921591258644
$ws.Range("L6").Value = 0.2968921591258644
$ws.Range("M6").Value = 0.299266677648297
$ws.Range("N6").Value = 0.299266677648297
$ws.Range("O6").Value = 0.5230334639612281
$ws.Range("P6").Value = 0.5230334639612281
$ws.Range("Q6").Value = 0.5230334639612281
$ws.Range("R6").Value = 0.6509088676058344
$ws.Range("S6").Value = 0.658957110924401
$ws.Range("T6").Value = 0.7938237279950101
$ws.Range("U6").Value = 0.8279793710301075
$ws.Range("V6").Value = 0.8574774278957599
$ws.Range("W6").Value = 0.8574774278957599
$ws.Range("X6").Value = 0.8574774278957599
$ws.Range("Y6").Value = 0.8574774278957599
$ws.Range("Z6").Value = 0.8692639309845254
$ws.Range("AA6").Value = 0.8692639309845254
$ws.Range("AB6").Value = 0.8692639309845254
$ws.Range("AC6").Value = 0.9514974023965904
$ws.Range("AD6").Value = 0.9715389119690973
$ws.Range("AE6").Value = 0.9999999999999997
$ws.Range("AF6").Value = 0.9999999999999997
$ws.Range("AG6").Value = 0.9999999999999997
$ws.Range("AH6").Value = 0.9999999999999997
$ws.Range("AI6").Value = 0.9999999999999997
$ws.Range("AJ6").Value = 0.9999999999999997

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5321267738290694
$ws.Range("F3").Value = 0.5550017396863993
$ws.Range("D4").Value = 13
$ws.Range("F4").Value = 0.5613680358433301
$ws.Range("G4").Value = 11
$ws.Range("F5").Value = 0.5701690915554083
$ws.Range("C6").Value = 3
$ws.Range("F6").Value = 0.5230334639612281
$ws.Range("G6").Value = 11

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.9148836870071987
$ws.Range("F3").Value = 0.8515397516793158
$ws.Range("F4").Value = 0.7303346303114059
$ws.Range("F5").Value = 0.8880190009284608
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 19
$ws.Range("F6").Value = 0.7938237279950101
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.9148836870071987
$ws.Range("F3").Value = 0.8515397516793158
$ws.Range("F4").Value = 0.9293446588695458
$ws.Range("F5").Value = 0.8880190009284608
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 20
$ws.Range("F6").Value = 0.8279793710301075
$ws.Range("G6").Value = 17

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9148836870071987
$ws.Range("D3").Value = 23
$ws.Range("F3").Value = 0.9009115337528599
$ws.Range("G3").Value = 22
$ws.Range("F4").Value = 0.9293446588695458
$ws.Range("D5").Value = 24
$ws.Range("F5").Value = 0.9075198200427399
$ws.Range("G5").Value = 22
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 28
$ws.Range("F6").Value = 0.9514974023965904
$ws.Range("G6").Value = 25

